# Auto-generated edit script applying cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "53.606.19"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -5.64%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.209.14"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -7.30%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "486.71"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "125.31"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.56%  "

$ws.Range("E7").Value = "  -0.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.522"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.237.87"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -6.62%  "

$ws.Range("E10").Value = "  -6.79%  "

$ws.Range("E11").Value = "  -0.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.319"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.64"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.604.80"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -7.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.25"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "53.547.54"
$ws.Range("D16").ClearFormats()

$ws.Range("E17").Value = "  -4.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.239.56"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -6.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.66"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.58%  "

$ws.Range("E20").Value = "  -2.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "295.77"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.15"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.71%  "

$ws.Range("E23").Value = "  -0.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.67"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.995"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("E26").Value = "  -1.23%  "

$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.312.81"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -7.38%  "

$ws.Range("E29").Value = "  -3.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "163.93"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.53%  "

$ws.Range("E31").Value = "  -4.55%  "

$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.79"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.31%  "

$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0670"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -6.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.993"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.66%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.30"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.46%  "

$ws.Range("E38").Value = "  -1.37%  "

$ws.Range("E39").Value = "  +0.70%  "

$ws.Range("E40").Value = "  -4.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.19"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.85%  "

$ws.Range("E42").Value = "  -1.17%  "

$ws.Range("E43").Value = "  -1.75%  "

$ws.Range("E44").Value = "  -3.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "126.65"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.82"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0881"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.536"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "233.39"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0471"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.54%  "

$ws.Range("E51").Value = "  -3.60%  "
